$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.936.07"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.536.04"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'317.12"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'96.45"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "'35.65"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'7.50"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("D14").Value = "2.927.32"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "2.499.79"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "'15.01"
$ws.Range("E16").Value = "  -6.08%  "
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").Value = "42.974.86"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").Value = "'12.54"
$ws.Range("E20").Value = "  -4.35%  "
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "'69.71"
$ws.Range("D23").Value = "'253.47"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").Value = "'27.03"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "'40.63"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").Value = "'10.32"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'5.88"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").Value = "'156.11"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'19.41"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").Value = "'3.36"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("D36").Value = "'0.0797"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("E38").Value = "  +1.35%  "
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  -7.59%  "
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "'0.0303"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("D46").Value = "2.001.13"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").Value = "'9.05"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Value = "'84.59"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "'74.98"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").Value = "2.782.19"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("E51").Value = "  +2.41%  "
